$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (was Anton / 1617.0 / Heptathlon 100 meters hurdles)
$ws.Range("A2").Value = "joel"
$ws.Range("B2").Value = 349.0
$ws.Range("C2").Value = "Decathlon 110 meters hurdles"

# Add new row 3
$ws.Range("A3").Value = "joel"
$ws.Range("B3").Value = 6941.0
$ws.Range("C3").Value = "Heptathlon Shot Put"

# Add new row 4
$ws.Range("A4").Value = "yes"
$ws.Range("B4").Value = 1606.0
$ws.Range("C4").Value = "Decathlon Long Jump"

# Add new row 5
$ws.Range("A5").Value = "hugo"
$ws.Range("B5").Value = 1237.0
$ws.Range("C5").Value = "Heptathlon High Jump"
